$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "electrode" column (D2:D6) from numeric 1 to text "A"
$ws.Range("D2").Value = "A"
$ws.Range("D3").Value = "A"
$ws.Range("D4").Value = "A"
$ws.Range("D5").Value = "A"
$ws.Range("D6").Value = "A"

# Update the selection shown in the sheet view
$ws.Range("J9:J10").Select()
